$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the degrees-of-freedom table (rows 2-6), replacing the previously
# blank/styled placeholder cells in A2:C4 with real data and extending it.
$ws.Range("A2:C4").ClearFormats()

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 1

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 2

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 4

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = 4

# Leftover styled (underlined) empty cell far below the table.
$ws.Range("D11").Font.Underline = $true

# Page setup for printing.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

$ws.Range("D6").Select()
